$d = $word.ActiveDocument

# The three list-item paragraphs (file paths) under "Archivos de la vista de
# reservas para gerencia" that must be highlighted in yellow:
#   -1)Routes/web.php
#   -2)Config/reservas.php
#   -3)app/Http/Controllers/GerenciaReservasController.php
# The following paragraph (-4)resources/views/.../index.blade.php) and the
# bold "Específicos para esta área:" divider must stay untouched.

$targets = @(
    "-1)Routes/web.php",
    "-2)Config/reservas.php",
    "-3)app/Http/Controllers/GerenciaReservasController.php"
)

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $text = $p.Range.Text
    foreach ($target in $targets) {
        if ($text -eq ($target + "`r")) {
            $r = $p.Range
            # Exclude the trailing paragraph-mark character so only the
            # visible run text (not the pilcrow) receives the highlight.
            [void]$r.MoveEnd(1, -1)
            $r.HighlightColorIndex = 7
        }
    }
}
